$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (before current row 290),
# shifting the existing rows 290-380 down to 292-382.
$ws.Range("A290:R291").Insert()

# Row 290: new weekly entry (Primera)
$ws.Cells.Item(290, 1).Value = 6
$ws.Cells.Item(290, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(290, 3).Value = "Metropolitana"
$ws.Cells.Item(290, 4).Value = "2021-09-24"
$ws.Cells.Item(290, 5).Value = 13
$ws.Cells.Item(290, 6).Value = 100112017
$ws.Cells.Item(290, 7).Value = "Apio"
$ws.Cells.Item(290, 8).Value = "Americana (o)"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 1400
$ws.Cells.Item(290, 11).Value = 7000
$ws.Cells.Item(290, 12).Value = 8000
$ws.Cells.Item(290, 13).Value = 7571
$ws.Cells.Item(290, 14).Value = "`$/docena de matas"
$ws.Cells.Item(290, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(290, 16).Value = 1262
$ws.Cells.Item(290, 17).Value = 6
$ws.Cells.Item(290, 18).Value = "Hortaliza"

# Row 291: new weekly entry (Segunda)
$ws.Cells.Item(291, 1).Value = 6
$ws.Cells.Item(291, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(291, 3).Value = "Metropolitana"
$ws.Cells.Item(291, 4).Value = "2021-09-24"
$ws.Cells.Item(291, 5).Value = 13
$ws.Cells.Item(291, 6).Value = 100112017
$ws.Cells.Item(291, 7).Value = "Apio"
$ws.Cells.Item(291, 8).Value = "Americana (o)"
$ws.Cells.Item(291, 9).Value = "Segunda"
$ws.Cells.Item(291, 10).Value = 600
$ws.Cells.Item(291, 11).Value = 6000
$ws.Cells.Item(291, 12).Value = 6000
$ws.Cells.Item(291, 13).Value = 6000
$ws.Cells.Item(291, 14).Value = "`$/docena de matas"
$ws.Cells.Item(291, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(291, 16).Value = 1000
$ws.Cells.Item(291, 17).Value = 6
$ws.Cells.Item(291, 18).Value = "Hortaliza"

Write-Host "done"
